# ===========================================================================
# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (fund holdings detail) right before the
#   "总计" (summary) worksheet.
# - Prepend a corresponding "2022-Q1" row to the "总计" summary sheet and
#   renumber the existing index column.
# ===========================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: Build the new "2022-Q1" detail sheet.
#
# The easiest way to get an identical header row / column styling to the
# other quarterly sheets is to duplicate the most recent quarter sheet
# ("2021-Q4") and then overwrite its data, placing the duplicate right
# before the "总计" sheet.
# ---------------------------------------------------------------------------

$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet.Copy($totalSheet)

# Re-fetch the freshly created copy by its (Excel-generated) name instead of
# relying on variables captured before the structural change, since sheet
# references can shift position once new sheets are inserted.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template sheet ("2021-Q4") has 20 data rows (rows 2-21). Clear any
# rows beyond that so only the 2022-Q1 rows remain.
$templateRows = $templateSheet.UsedRange.Rows.Count
for ($r = 22; $r -le $templateRows; $r++) {
    $newSheet.Rows.Item($r).Clear()
}

# The fund-code / numeric-text columns (B and D:G) must stay plain text so
# values such as "000452" keep their leading zeros and values such as
# "31.88" are not silently turned into floating point numbers.
$newSheet.Range("B2:B21").NumberFormat = "@"
$newSheet.Range("D2:G21").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "000452"
$newSheet.Range("C2").Value = "南方医药保健灵活配置混合"
$newSheet.Range("D2").Value = "31.88"
$newSheet.Range("E2").Value = "90.98"
$newSheet.Range("F2").Value = "3.52"
$newSheet.Range("G2").Value = "1.1222"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "010592"
$newSheet.Range("C3").Value = "南方医药创新股票A"
$newSheet.Range("D3").Value = "23.21"
$newSheet.Range("E3").Value = "87.52"
$newSheet.Range("F3").Value = "3.43"
$newSheet.Range("G3").Value = "0.7961"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "008934"
$newSheet.Range("C4").Value = "大成科技消费股票A"
$newSheet.Range("D4").Value = "13.43"
$newSheet.Range("E4").Value = "83.81"
$newSheet.Range("F4").Value = "3.46"
$newSheet.Range("G4").Value = "0.4647"
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "003230"
$newSheet.Range("C5").Value = "创金合信医疗保健行业股票A"
$newSheet.Range("D5").Value = "8.07"
$newSheet.Range("E5").Value = "94.55"
$newSheet.Range("F5").Value = "4.61"
$newSheet.Range("G5").Value = "0.3720"
$newSheet.Range("H5").Value = 9

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "001766"
$newSheet.Range("C6").Value = "上投摩根医疗健康股票"
$newSheet.Range("D6").Value = "10.35"
$newSheet.Range("E6").Value = "80.54"
$newSheet.Range("F6").Value = "3.31"
$newSheet.Range("G6").Value = "0.3426"
$newSheet.Range("H6").Value = 6

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "010593"
$newSheet.Range("C7").Value = "南方医药创新股票C"
$newSheet.Range("D7").Value = "8.05"
$newSheet.Range("E7").Value = "87.52"
$newSheet.Range("F7").Value = "3.43"
$newSheet.Range("G7").Value = "0.2761"
$newSheet.Range("H7").Value = 10

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "003231"
$newSheet.Range("C8").Value = "创金合信医疗保健行业股票C"
$newSheet.Range("D8").Value = "4.28"
$newSheet.Range("E8").Value = "94.55"
$newSheet.Range("F8").Value = "4.61"
$newSheet.Range("G8").Value = "0.1973"
$newSheet.Range("H8").Value = 9

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "010585"
$newSheet.Range("C9").Value = "创金合信医药消费股票A"
$newSheet.Range("D9").Value = "4.22"
$newSheet.Range("E9").Value = "93.65"
$newSheet.Range("F9").Value = "4.42"
$newSheet.Range("G9").Value = "0.1865"
$newSheet.Range("H9").Value = 9

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "161039"
$newSheet.Range("C10").Value = "富国中证1000指数增强LOF"
$newSheet.Range("D10").Value = "21.72"
$newSheet.Range("E10").Value = "89.03"
$newSheet.Range("F10").Value = "0.65"
$newSheet.Range("G10").Value = "0.1412"
$newSheet.Range("H10").Value = 10

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "008935"
$newSheet.Range("C11").Value = "大成科技消费股票C"
$newSheet.Range("D11").Value = "3.98"
$newSheet.Range("E11").Value = "83.81"
$newSheet.Range("F11").Value = "3.46"
$newSheet.Range("G11").Value = "0.1377"
$newSheet.Range("H11").Value = 9

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "002210"
$newSheet.Range("C12").Value = "创金合信量化多因子股票A"
$newSheet.Range("D12").Value = "3.31"
$newSheet.Range("E12").Value = "88.74"
$newSheet.Range("F12").Value = "0.89"
$newSheet.Range("G12").Value = "0.0295"
$newSheet.Range("H12").Value = 10

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "003594"
$newSheet.Range("C13").Value = "长盛盛崇灵活配置混合A"
$newSheet.Range("D13").Value = "1.84"
$newSheet.Range("E13").Value = "43.63"
$newSheet.Range("F13").Value = "1.42"
$newSheet.Range("G13").Value = "0.0261"
$newSheet.Range("H13").Value = 10

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "161038"
$newSheet.Range("C14").Value = "富国新兴成长量化精选混合（LOF）"
$newSheet.Range("D14").Value = "1.13"
$newSheet.Range("E14").Value = "93.66"
$newSheet.Range("F14").Value = "2.01"
$newSheet.Range("G14").Value = "0.0227"
$newSheet.Range("H14").Value = 5

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "010586"
$newSheet.Range("C15").Value = "创金合信医药消费股票C"
$newSheet.Range("D15").Value = "0.50"
$newSheet.Range("E15").Value = "93.65"
$newSheet.Range("F15").Value = "4.42"
$newSheet.Range("G15").Value = "0.0221"
$newSheet.Range("H15").Value = 9

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "080008"
$newSheet.Range("C16").Value = "长盛战略新兴产业灵活配置混合A"
$newSheet.Range("D16").Value = "1.83"
$newSheet.Range("E16").Value = "40.84"
$newSheet.Range("F16").Value = "1.19"
$newSheet.Range("G16").Value = "0.0218"
$newSheet.Range("H16").Value = 10

$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "001834"
$newSheet.Range("C17").Value = "长盛战略新兴产业灵活配置混合C"
$newSheet.Range("D17").Value = "1.43"
$newSheet.Range("E17").Value = "40.84"
$newSheet.Range("F17").Value = "1.19"
$newSheet.Range("G17").Value = "0.0170"
$newSheet.Range("H17").Value = 10

$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "014285"
$newSheet.Range("C18").Value = "鑫元健康产业混合A"
$newSheet.Range("D18").Value = "0.40"
$newSheet.Range("E18").Value = "33.76"
$newSheet.Range("F18").Value = "2.22"
$newSheet.Range("G18").Value = "0.0089"
$newSheet.Range("H18").Value = 5

$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "003865"
$newSheet.Range("C19").Value = "创金合信量化多因子股票C"
$newSheet.Range("D19").Value = "0.79"
$newSheet.Range("E19").Value = "88.74"
$newSheet.Range("F19").Value = "0.89"
$newSheet.Range("G19").Value = "0.0070"
$newSheet.Range("H19").Value = 10

$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "014286"
$newSheet.Range("C20").Value = "鑫元健康产业混合C"
$newSheet.Range("D20").Value = "0.27"
$newSheet.Range("E20").Value = "33.76"
$newSheet.Range("F20").Value = "2.22"
$newSheet.Range("G20").Value = "0.0060"
$newSheet.Range("H20").Value = 5

$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "003595"
$newSheet.Range("C21").Value = "长盛盛崇灵活配置混合C"
$newSheet.Range("D21").Value = "0.11"
$newSheet.Range("E21").Value = "43.63"
$newSheet.Range("F21").Value = "1.42"
$newSheet.Range("G21").Value = "0.0016"
$newSheet.Range("H21").Value = 10

# ---------------------------------------------------------------------------
# Step 2: Update the "总计" (summary) sheet - insert a new row for 2022-Q1
# above the existing "2021-Q4" row, shifting the rest down, and renumber the
# index column (A) so it continues to read 0,1,2,3,4,5.
# ---------------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Copy the formatting of the (still untouched) index cell in row 3 onto the
# freshly inserted row 2 index cell, so it keeps the same bold/border style
# used throughout column A, then clear any stray formatting the insert
# operation may have copied into B2:D2 from the row below.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 20
$totalSheet.Range("D2").Value = 4.2

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
